$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.334.99"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "2.283.19"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.61"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.64"
$ws.Range("E6").Value = "  -1.29%  "

$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.77"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.01"
$ws.Range("E12").Value = "  +3.74%  "

$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("E14").Value = "  -0.89%  "

$ws.Range("D15").Value = "2.623.35"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.866"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").Value = "2.281.67"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").Value = "43.209.67"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.84"
$ws.Range("E20").Value = "  +3.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.56"
$ws.Range("E21").Value = "  -1.76%  "

$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.34"
$ws.Range("E23").Value = "  -1.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.64"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("E27").Value = "  -1.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.88"
$ws.Range("E28").Value = "  -8.04%  "

$ws.Range("E29").Value = "  -2.12%  "

$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.24"
$ws.Range("E31").Value = "  -3.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.44"
$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0907"
$ws.Range("E33").Value = "  -3.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.78"
$ws.Range("E34").Value = "  +4.02%  "

$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.93"

$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("E39").Value = "  -4.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  +10.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.07"
$ws.Range("E41").Value = "  +3.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.95"
$ws.Range("E42").Value = "  +3.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("E44").Value = "  +4.49%  "

$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.71"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.14"
$ws.Range("E48").Value = "  +2.01%  "

$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0997"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.438"
$ws.Range("E51").Value = "  -3.79%  "
